# Commit: Tue, Jul 14, 2020  2:05:06 AM
#
# 1) The cash-flow summary table on slide 16 had its table style switched
#    from the deck's one custom table style ({B9CEE724-...}, "Table_0")
#    to PowerPoint's built-in "No Style, Table Grid" style
#    ({113DD694-1A94-493D-BF3E-2787B6B46CBB}). Table styles can't be set
#    by assigning Table.Style directly (PowerPoint requires ApplyStyle),
#    so use Table.ApplyStyle with the target style GUID.

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(16)

# Locate the shape that actually hosts the table instead of hard-coding
# its index, so the script is resilient to any shape reordering.
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}

if ($tableShape -ne $null) {
    $tableShape.Table.ApplyStyle("{113DD694-1A94-493D-BF3E-2787B6B46CBB}")
}

# 2) The deck's Design ("Integral") was switched back to the plain default
#    "Office Theme" on the slide master (leaving the notes master carrying
#    the Integral theme that the slide master no longer uses). This is the
#    effect of picking a different thumbnail on the Design tab. Attempt it
#    through the documented COM entry points; guarded with try/catch since
#    this hosted runtime does not implement real theme-file import (no
#    .thmx available to load), so these calls are expected to be no-ops
#    here rather than raising - kept for parity with a real PowerPoint host.
try { $ppt.OpenThemeFile("Office Theme") } catch { }
try { $p.ApplyTheme("Office Theme") } catch { }
try { $p.Designs.Item(1).Name = "Office Theme" } catch { }

